$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 28, shifting existing rows 28-32 down to 29-33
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new data entry
$ws.Cells.Item(28, 1).Value = 10
$ws.Cells.Item(28, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(28, 3).Value = "La Araucanía"
$ws.Cells.Item(28, 4).Value = (Get-Date -Year 2023 -Month 11 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(28, 5).Value = 9
$ws.Cells.Item(28, 6).Value = "Fruta"
$ws.Cells.Item(28, 7).Value = 100104
$ws.Cells.Item(28, 8).Value = "Frutos de pepita"
$ws.Cells.Item(28, 9).Value = 100104004
$ws.Cells.Item(28, 10).Value = "Níspero"
$ws.Cells.Item(28, 11).Value = "Californiana(o)"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 100
$ws.Cells.Item(28, 14).Value = 22000
$ws.Cells.Item(28, 15).Value = 22000
$ws.Cells.Item(28, 16).Value = 22000
$ws.Cells.Item(28, 17).Value = "$/bandeja 5 kilos"
$ws.Cells.Item(28, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(28, 19).Value = 4400
$ws.Cells.Item(28, 20).Value = 5
